$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty placeholder inline-string cells in column F (rows 2-17).
# F12 already holds real content ("ADDICTO: ...") and must stay untouched.
$ws.Range("F2:F11").ClearContents()
$ws.Range("F13:F17").ClearContents()

# Append a new row (18) importing "personal attribute [BCIO:050300]"
$ws.Range("A18").Value = "bcio"
$ws.Range("C18").Value = "entity [BFO:0000001]"
$ws.Range("D18").Value = "personal attribute [BCIO:050300]"
$ws.Range("E18").Value = "all"
